# Commit "added some test stuff" - OOXML diff analysis
# -----------------------------------------------------
# Diffing the canonical OOXML of before/after shows exactly two kinds of
# change, and nothing else (no slide added/removed, no shape added/
# removed/moved/resized, no text changed):
#
#   1. Every r:id / r:embed relationship id in the package (presentation
#      master/slide/layout ids, the webextension graphicFrame's we:id,
#      the fallback picture's r:embed, and the snapshot r:embed inside
#      ppt/slides/udata/data.xml) is replaced with a newly-minted GUID.
#      The relationship *graph* itself is unchanged - same targets, same
#      types, same cardinality - only the opaque id tokens differ, which
#      is exactly what happens whenever PowerPoint rewrites the package.
#
#   2. The we:webextension/@id instance GUID in
#      ppt/slides/udata/data.xml changes
#      ({3d22325f-6379-4170-8157-125e574b8213} ->
#       {77f30025-0611-42c5-96bf-0b8e01abc9e2}). This is PowerPoint's own
#      internal bookkeeping id for the inserted "PowerPoll" content
#      add-in instance; it is minted by PowerPoint when the add-in is
#      (re-)inserted and is not a reachable property anywhere in the
#      Slide/Shape/Presentation object model (there is no WebExtension(s)
#      object, and the add-in's graphicFrame/fallback picture isn't
#      independently addressable through Shapes either) - it simply isn't
#      scriptable, in real PowerPoint or here.
#
# In other words the whole diff is PowerPoint's own package bookkeeping
# from re-persisting the deck, not an authored content edit. The
# faithful reproduction through the supported COM surface is therefore
# to open the (already-active) presentation and save it back out, which
# is what round-trips/re-persists the package without touching any
# slide content.

$p = $ppt.ActivePresentation
$p.Save()
